# Auto-generated Excel COM-interop script
# Applies the scheduled-runner price/profit refresh to the Anima_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2833.3333
$ws.Range("J51").Value = 2833.3333
$ws.Range("L51").Value = 2833.3333
$ws.Range("N51").Value = -3801.3333
# Row 62
$ws.Range("H62").Value = 2290.25
$ws.Range("I62").Value = 1354.8572
$ws.Range("J62").Value = 3599.8
$ws.Range("K62").Value = 1354.8572
$ws.Range("L62").Value = 3599.8
$ws.Range("M62").Value = -730.8571999999999
$ws.Range("N62").Value = -4847.8
# Row 65
$ws.Range("H65").Value = 2290.25
$ws.Range("I65").Value = 1354.8572
$ws.Range("J65").Value = 3599.8
$ws.Range("K65").Value = 6774.286
$ws.Range("L65").Value = 17999
$ws.Range("M65").Value = -3654.286
$ws.Range("N65").Value = -24239
# Row 121
$ws.Range("H121").Value = 1698.5454
$ws.Range("J121").Value = 1769.5807
$ws.Range("L121").Value = 5308.742099999999
$ws.Range("N121").Value = -8802.742099999999
# Row 138
$ws.Range("H138").Value = 3317.6287
$ws.Range("I138").Value = 3633
$ws.Range("J138").Value = 3173.0833
$ws.Range("K138").Value = 10899
$ws.Range("L138").Value = 9519.249899999999
$ws.Range("M138").Value = -5759
$ws.Range("N138").Value = -19799.2499

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 934.3
$ws.Range("I2").Value = 809.46155
$ws.Range("J2").Value = 1166.1428
$ws.Range("K2").Value = 809.46155
$ws.Range("L2").Value = 1166.1428
$ws.Range("M2").Value = -696.46155
$ws.Range("N2").Value = -1392.1428
# Row 32
$ws.Range("H32").Value = 337622.34
$ws.Range("I32").Value = 388660.7
$ws.Range("K32").Value = 388660.7
$ws.Range("M32").Value = -388373.7
# Row 116
$ws.Range("H116").Value = 934.3
$ws.Range("I116").Value = 809.46155
$ws.Range("J116").Value = 1166.1428
$ws.Range("K116").Value = 809.46155
$ws.Range("L116").Value = 1166.1428
$ws.Range("M116").Value = 1484.53845
$ws.Range("N116").Value = -5754.1428
# Row 122
$ws.Range("H122").Value = 1460.4615
$ws.Range("I122").Value = 1415.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4246.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1796.5
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 934.3
$ws.Range("I3").Value = 809.46155
$ws.Range("J3").Value = 1166.1428
$ws.Range("K3").Value = 809.46155
$ws.Range("L3").Value = 1166.1428
$ws.Range("M3").Value = -695.46155
$ws.Range("N3").Value = -1394.1428
# Row 86
$ws.Range("H86").Value = 1537.6666
$ws.Range("I86").Value = 1457
$ws.Range("J86").Value = 1699
$ws.Range("K86").Value = 1457
$ws.Range("L86").Value = 1699
$ws.Range("M86").Value = -334
$ws.Range("N86").Value = -3945
# Row 89
$ws.Range("H89").Value = 1537.6666
$ws.Range("I89").Value = 1457
$ws.Range("J89").Value = 1699
$ws.Range("K89").Value = 7285
$ws.Range("L89").Value = 8495
$ws.Range("M89").Value = -1669
$ws.Range("N89").Value = -19727
# Row 110
$ws.Range("H110").Value = 100702
$ws.Range("J110").Value = 100702
$ws.Range("L110").Value = 100702
$ws.Range("N110").Value = -108882

$ws = $wb.Worksheets.Item("CRP")
# Row 44
$ws.Range("H44").Value = 18035.5
$ws.Range("J44").Value = 18035.5
$ws.Range("L44").Value = 18035.5
$ws.Range("N44").Value = -18919.5
# Row 86
$ws.Range("H86").Value = 1728.7368
$ws.Range("I86").Value = 1848.28
$ws.Range("K86").Value = 1848.28
$ws.Range("M86").Value = -725.28
# Row 89
$ws.Range("H89").Value = 1728.7368
$ws.Range("I89").Value = 1848.28
$ws.Range("K89").Value = 9241.4
$ws.Range("M89").Value = -3625.4

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1245.25
$ws.Range("I68").Value = 1005.8333
$ws.Range("J68").Value = 1364.9584
$ws.Range("K68").Value = 3017.4999
$ws.Range("L68").Value = 4094.8752
$ws.Range("M68").Value = -2206.4999
$ws.Range("N68").Value = -5716.8752
# Row 71
$ws.Range("H71").Value = 1245.25
$ws.Range("I71").Value = 1005.8333
$ws.Range("J71").Value = 1364.9584
$ws.Range("K71").Value = 9052.4997
$ws.Range("L71").Value = 12284.6256
$ws.Range("M71").Value = -4996.4997
$ws.Range("N71").Value = -20396.6256
# Row 114
$ws.Range("H114").Value = 905.625
$ws.Range("I114").Value = 186.5
$ws.Range("J114").Value = 1624.75
$ws.Range("K114").Value = 559.5
$ws.Range("L114").Value = 4874.25
$ws.Range("M114").Value = 2694.5
$ws.Range("N114").Value = -11382.25
# Row 122
$ws.Range("H122").Value = 5662.1
$ws.Range("I122").Value = 335.9091
$ws.Range("J122").Value = 12171.889
$ws.Range("K122").Value = 3023.1819
$ws.Range("L122").Value = 109547.001
$ws.Range("M122").Value = -573.1819
$ws.Range("N122").Value = -114447.001
# Row 129
$ws.Range("H129").Value = 641.3333
$ws.Range("I129").Value = 404.875
$ws.Range("J129").Value = 2533
$ws.Range("K129").Value = 1214.625
$ws.Range("L129").Value = 7599
$ws.Range("M129").Value = 3785.375
$ws.Range("N129").Value = -17599
# Row 131
$ws.Range("H131").Value = 1151.125
$ws.Range("I131").Value = 1038.9231
$ws.Range("J131").Value = 1227.8948
$ws.Range("K131").Value = 3116.7693
$ws.Range("L131").Value = 3683.6844
$ws.Range("M131").Value = 1923.2307
$ws.Range("N131").Value = -13763.6844
# Row 136
$ws.Range("H136").Value = 3252.5334
$ws.Range("I136").Value = 3028.8
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 9086.400000000001
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -3986.400000000001
$ws.Range("N136").Value = -21300
# Row 137
$ws.Range("H137").Value = 8771.058999999999
$ws.Range("I137").Value = 12952.8
$ws.Range("J137").Value = 2797.1428
$ws.Range("K137").Value = 38858.39999999999
$ws.Range("L137").Value = 8391.428400000001
$ws.Range("M137").Value = -33758.39999999999
$ws.Range("N137").Value = -18591.4284
# Row 141
$ws.Range("H141").Value = 21751.5
$ws.Range("I141").Value = 15006
$ws.Range("J141").Value = 32994
$ws.Range("K141").Value = 45018
$ws.Range("L141").Value = 98982
$ws.Range("M141").Value = -39838
$ws.Range("N141").Value = -109342

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1529.8889
$ws.Range("I102").Value = 1529.8889
$ws.Range("K102").Value = 1529.8889
$ws.Range("M102").Value = 92.11110000000008
# Row 132
$ws.Range("H132").Value = 2636.182
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -17059.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 166669500
$ws.Range("I7").Value = 200002400
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 200002400
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -200002288
$ws.Range("N7").Value = -5229
# Row 22
$ws.Range("H22").Value = 9798.916999999999
$ws.Range("I22").Value = 696
$ws.Range("J22").Value = 16301
$ws.Range("K22").Value = 696
$ws.Range("L22").Value = 16301
$ws.Range("M22").Value = -401
$ws.Range("N22").Value = -16891
# Row 27
$ws.Range("H27").Value = 9798.916999999999
$ws.Range("I27").Value = 696
$ws.Range("J27").Value = 16301
$ws.Range("K27").Value = 696
$ws.Range("L27").Value = 16301
$ws.Range("M27").Value = -589
$ws.Range("N27").Value = -16515
# Row 121
$ws.Range("H121").Value = 80308.89
$ws.Range("J121").Value = 80308.89
$ws.Range("L121").Value = 80308.89
$ws.Range("N121").Value = -83802.89
# Row 126
$ws.Range("H126").Value = 166669500
$ws.Range("I126").Value = 200002400
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 600007200
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -600004730
$ws.Range("N126").Value = -19955
# Row 132
$ws.Range("H132").Value = 3234.0442
$ws.Range("I132").Value = 3331.9092
$ws.Range("J132").Value = 3141.7715
$ws.Range("K132").Value = 9995.7276
$ws.Range("L132").Value = 9425.3145
$ws.Range("M132").Value = -7465.7276
$ws.Range("N132").Value = -14485.3145

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4312.231
$ws.Range("I81").Value = 4395.5557
$ws.Range("K81").Value = 8791.1114
$ws.Range("M81").Value = -7730.1114
# Row 84
$ws.Range("H84").Value = 4312.231
$ws.Range("I84").Value = 4395.5557
$ws.Range("K84").Value = 43955.557
$ws.Range("M84").Value = -38651.557
# Row 103
$ws.Range("H103").Value = 70000
$ws.Range("J103").Value = 70000
$ws.Range("L103").Value = 70000
$ws.Range("N103").Value = -72344
# Row 112
$ws.Range("H112").Value = 60693.5
$ws.Range("J112").Value = 60693.5
$ws.Range("L112").Value = 60693.5
$ws.Range("N112").Value = -63647.5
# Row 121
$ws.Range("H121").Value = 50420
$ws.Range("J121").Value = 50420
$ws.Range("L121").Value = 50420
$ws.Range("N121").Value = -53914
# Row 122
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3125
$ws.Range("K122").Value = 9375
$ws.Range("M122").Value = -6925
# Row 126
$ws.Range("H126").Value = 1820.3636
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1820.3636
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5461.0908
$ws.Range("N126").Value = -10401.0908
$ws.Range("M126").ClearContents()
# Row 136
$ws.Range("H136").Value = 2085
$ws.Range("I136").Value = 1838.0851
$ws.Range("J136").Value = 2665.25
$ws.Range("K136").Value = 5514.2553
$ws.Range("L136").Value = 7995.75
$ws.Range("M136").Value = -2964.2553
$ws.Range("N136").Value = -13095.75
